$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '58.372.73'
$ws.Range('E2').Value2 = '  -3.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.695.77'
$ws.Range('E3').Value2 = '  -7.05%  '
$ws.Range('E4').Value2 = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '499.86'
$ws.Range('E5').Value2 = '  -5.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '139.05'
$ws.Range('E6').Value2 = '  -3.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.998'
$ws.Range('E7').Value2 = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.525'
$ws.Range('E8').Value2 = '  -5.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '2.707.20'
$ws.Range('E9').Value2 = '  -6.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '6.03'
$ws.Range('E10').Value2 = '  +1.23%  '
$ws.Range('E11').Value2 = '  -3.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.346'
$ws.Range('E12').Value2 = '  -4.62%  '
$ws.Range('E13').Value2 = '  +0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '3.168.30'
$ws.Range('E14').Value2 = '  -7.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '58.435.48'
$ws.Range('E15').Value2 = '  -3.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '21.37'
$ws.Range('E16').Value2 = '  -5.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '2.703.86'
$ws.Range('E17').Value2 = '  -6.75%  '
$ws.Range('E18').Value2 = '  -6.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '4.71'
$ws.Range('E19').Value2 = '  -6.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '10.91'
$ws.Range('E20').Value2 = '  -7.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '333.26'
$ws.Range('E21').Value2 = '  -8.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '6.19'
$ws.Range('E22').Value2 = '  -6.99%  '
$ws.Range('E23').Value2 = '  -0.41%  '
$ws.Range('E24').Value2 = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '62.77'
$ws.Range('E25').Value2 = '  -2.49%  '
$ws.Range('B26').Value2 = 'Polygon'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '0.423'
$ws.Range('E26').Value2 = '  -6.74%  '
$ws.Range('B27').Value2 = 'Kaspa'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '0.171'
$ws.Range('E27').Value2 = '  -4.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '0.997'
$ws.Range('E28').Value2 = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '7.38'
$ws.Range('E29').Value2 = '  -5.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '0.0₃0815'
$ws.Range('E30').Value2 = '  -5.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '0.999'
$ws.Range('E31').Value2 = '  -0.02%  '
$ws.Range('B32').Value2 = 'EthereumClassic'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '19.05'
$ws.Range('E32').Value2 = '  -3.39%  '
$ws.Range('B33').Value2 = 'PancakeSwap'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '1.59'
$ws.Range('E33').Value2 = '  -5.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '150.55'
$ws.Range('E34').Value2 = '  +1.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '5.35'
$ws.Range('E35').Value2 = '  -4.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '4.13'
$ws.Range('E36').Value2 = '  -5.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '0.927'
$ws.Range('E37').Value2 = '  -7.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '1.11'
$ws.Range('E38').Value2 = '  -8.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '35.32'
$ws.Range('E39').Value2 = '  -6.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '3.54'
$ws.Range('E40').Value2 = '  -4.19%  '
$ws.Range('B41').Value2 = 'Maker'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '2.174.04'
$ws.Range('E41').Value2 = '  -6.45%  '
$ws.Range('B42').Value2 = 'Stacks'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '1.36'
$ws.Range('E42').Value2 = '  -9.12%  '
$ws.Range('E43').Value2 = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.0555'
$ws.Range('E44').Value2 = '  -4.53%  '
$ws.Range('E45').Value2 = '  -8.02%  '
$ws.Range('E46').Value2 = '  +0.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '18.63'
$ws.Range('E47').Value2 = '  -10.07%  '
$ws.Range('B48').Value2 = 'VeChain'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.0225'
$ws.Range('E48').Value2 = '  -4.53%  '
$ws.Range('B49').Value2 = 'RenderToken'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '4.60'
$ws.Range('E49').Value2 = '  -7.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.0883'
$ws.Range('E50').Value2 = '  -5.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '17.73'
$ws.Range('E51').Value2 = '  -4.39%  '
